# The workbook was re-saved from a Korean-locale copy of Excel (South Korea
# localization pass). The only actual content change in that re-save is the
# footnote on the "About" sheet: the US-specific note is swapped for the
# Korea-specific one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("A9").Value = "In Korea, we set this value to 0."
